# Phase-3 budget workbook cleanup:
#  - drop the now-unused "api calls" and "Sheet3" worksheets (and the
#    shared strings that were only referenced there)
#  - shift the trailing "Total hours spent so far" summary row down two
#    rows (36 -> 38), leaving two blank rows above it
#  - update the active selection to reflect where editing left off

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("api calls").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("begroting")

# Push the final totals row (previously row 36) down to row 38, inserting
# two blank rows above it; formulas referencing it are auto-adjusted.
$ws.Rows("36:37").Insert()

$ws.Range("C30").Select()
